# Learners workbook update:
#  - replace placeholder learner rows (Sheldon Cooper / Ted Mosby) with
#    real learner data (Akshay / Vinod)
#  - remove the Text number format previously forced on column C (Phone Number)
#    so the phone numbers are stored as plain numbers
#  - add real "mailto:" hyperlinks on the Email column cells
#  - move the active selection to E3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C ("Phone Number") was formatted as Text (numFmt 49). Clear that
# formatting (on the whole column and the already-used cells) before writing
# the new numeric phone numbers so they are stored as numbers, not strings.
$ws.Columns.Item(3).ClearFormats()

# Row 2: Akshay / IT / phone / email
$ws.Range("A2").Value = "Akshay"
$ws.Range("B2").Value = "IT"
$ws.Range("C2").Value = 7656787890
$ws.Range("D2").Value = "akshay@gmail.com"

# Row 3: Vinod / HR / phone / email
$ws.Range("A3").Value = "Vinod"
$ws.Range("B3").Value = "HR"
$ws.Range("C3").Value = 7869352434
$ws.Range("D3").Value = "vinod@gmail.com"

# Turn the email addresses into real hyperlinks (mailto:)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:akshay@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:vinod@gmail.com")

# Leave the selection on E3, matching the saved view state
$ws.Range("E3").Select() | Out-Null
